$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'basketball gear boys'
$ws.Range("A2").Value = 'spandex shorts for volleyball'
$ws.Range("A3").Value = 'volleyball knee sleeves men'
$ws.Range("A4").Value = 'knee pads for work construction'
$ws.Range("A5").Value = 'knee pads for working on floors'
$ws.Range("A6").Value = 'bee tights adult'
$ws.Range("A7").Value = 'compression test'
$ws.Range("A8").Value = 'lacrosse youth shorts'
$ws.Range("A9").Value = 'women volleyball knee pads'
$ws.Range("A10").Value = 'short football pants'
$ws.Range("A11").Value = 'girls compression tights'
$ws.Range("A12").Value = 'knee compression sleeve reduce strain & swelling'
$ws.Range("A13").Value = 'volleyball knee pads small'
$ws.Range("A14").Value = 'yoga pants knee length'
$ws.Range("A15").Value = 'pad for squats'
$ws.Range("A16").Value = 'sweat pads'
$ws.Range("A17").Value = 'calf protector'
$ws.Range("A18").Value = 'compression pants long'
$ws.Range("A19").Value = 'spandex mens pants'
$ws.Range("A20").Value = 'protective knee pads for work'
$ws.Range("A21").Value = 'youth compression shorts'
$ws.Range("A22").Value = 'mens basketball shorts long'
$ws.Range("A23").Value = 'thick volleyball knee pads'
$ws.Range("A24").Value = 'weightlifting guide'
$ws.Range("A25").Value = 'small work knee pads'
$ws.Range("A26").Value = 'basketballs in bulk'
$ws.Range("A27").Value = 'baseballs cheap'
$ws.Range("A28").Value = 'spandex capri leggings'
$ws.Range("A29").Value = 'hex fabric'
$ws.Range("A30").Value = 'black knee pads for work'
$ws.Range("A31").Value = 'black knee guards'
$ws.Range("A32").Value = 'work knee pads for men gel'
$ws.Range("A33").Value = 'leg guard baseball'
$ws.Range("A34").Value = 'lacrosse shorts youth'
$ws.Range("A35").Value = 'girdles for men'
$ws.Range("A36").Value = 'basketball cheap'
$ws.Range("A37").Value = 'wrestling fight shorts'
$ws.Range("A38").Value = 'construction knee'
$ws.Range("A39").Value = 'knee compression sleeve protector'
$ws.Range("A40").Value = 'knee bursitis sleeve'
$ws.Range("A41").Value = 'football leggings'
$ws.Range("A42").Value = 'black girls softball pants'
$ws.Range("A43").Value = 'girls softball pants black'
$ws.Range("A44").Value = 'working pants with knee pads'
$ws.Range("A45").Value = 'youth volleyball sleeves'
$ws.Range("A46").Value = 'mens knee length shorts'
$ws.Range("A47").Value = 'compression capris girls'
$ws.Range("A48").Value = 'leaf leggings'
$ws.Range("A49").Value = 'fit compression knee'
$ws.Range("A50").Value = 'performance basketball'
$ws.Range("A51").Value = 'pantalones de basketball'
$ws.Range("A52").Value = 'youth baseball sleeves for boys'
$ws.Range("A53").Value = 'mens tights for sports'
$ws.Range("A54").Value = 'knee protector work'
$ws.Range("A55").Value = 'long basketball shorts for men'
$ws.Range("A56").Value = 'mens knee sleeves weightlifting'
$ws.Range("A57").Value = 'knee pads bulk'
$ws.Range("A58").Value = 'snowboarding protective gear'
$ws.Range("A59").Value = 'calf compression sleeve youth'
$ws.Range("A60").Value = 'knee pain pads'
$ws.Range("A61").Value = 'knee pad for construction'
$ws.Range("A62").Value = 'mens capri yoga pants'
$ws.Range("A63").Value = 'knee sleeves for basketball'
$ws.Range("A64").Value = 'knee pad sleeves'
$ws.Range("A65").Value = 'knee work'
$ws.Range("A66").Value = 'volleyball shorts men'
$ws.Range("A67").Value = 'knee pad for working'
$ws.Range("A68").Value = 'above knee shorts men'
$ws.Range("A69").Value = 'youth volleyball shorts for girls'
$ws.Range("A70").Value = 'protector paintball'
$ws.Range("A71").Value = 'baseball shorts boys'
$ws.Range("A72").Value = 'boys tights and leggings'
$ws.Range("A73").Value = 'kneeling pad gel'
$ws.Range("A74").Value = 'paintball pants men'
$ws.Range("A75").Value = 'men above knee shorts'
$ws.Range("A76").Value = 'spandex compression shorts'
$ws.Range("A77").Value = 'hockey padded shorts'
$ws.Range("A78").Value = 'knee pads for joint pain'
$ws.Range("A79").Value = 'rodillera volleyball'
$ws.Range("A80").Value = 'basketball calf sleeve'
$ws.Range("A81").Value = 'womens lacrosse pants'
$ws.Range("A82").Value = 'best work knee pads'
$ws.Range("A83").Value = 'leg sleeves for men basketball'
$ws.Range("A84").Value = 'football pouch youth'
$ws.Range("A85").Value = 'calf compression sleeve boys'
$ws.Range("A86").Value = 'boys leggings sports'
$ws.Range("A87").Value = 'basketball stretch pants'
$ws.Range("A88").Value = 'mens spandex shorts'
$ws.Range("A89").Value = 'black kneepads'
$ws.Range("A90").Value = 'joint protectors'
$ws.Range("A91").Value = 'boys hiking pants'
$ws.Range("A92").Value = 'knee length yoga pants'
$ws.Range("A93").Value = 'girls volleyball shorts youth'
$ws.Range("A94").Value = 'baseball pants mens long'
$ws.Range("A95").Value = 'mens sports tights'
$ws.Range("A96").Value = 'womens football pads'
$ws.Range("A97").Value = 'stretch mark men'
$ws.Range("A98").Value = 'youth xl baseball pants'
$ws.Range("A99").Value = 'youth basketball compression sleeve'
$ws.Range("A100").Value = 'mens gym tights'
